$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Data Inizio (D) and Data Fine (E) columns for rows 2-6 to the new date/time values
$dataInizio = 45900.272222222222
$dataFine = 45900.51666666667

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 4).Value = $dataInizio
    $ws.Cells.Item($r, 5).Value = $dataFine
}

# Update the active selection from D5 to D4
$ws.Range("D4").Select()
